$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),
                (''model'',
                 LogisticRegression(l1_ratio=0.7, max_iter=1000,
                                    penalty=''elasticnet'', random_state=42,
                                    solver=''saga''))])'
$ws.Range("B2").Value = 0.675
$ws.Range("C2").Value = '{''selector'': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), ''scaler'': MinMaxScaler(), ''model__solver'': ''saga'', ''model__penalty'': ''elasticnet'', ''model__l1_ratio'': 0.7, ''model__class_weight'': None}'
$ws.Range("D2").Value = 0.75
$ws.Range("E2").Value = '[1 0 1 1 1 1 0 1 0 1 0 1]'
$ws.Range("F2").Value = '[1 1 1 1 1 1 0 0 1 0 0 1]'
$ws.Range("G2").Value = 42
$ws.Range("H2").Value = 0.7226190422428765
$ws.Range("I2").Value = 0.02060886693193312
$ws.Range("J2").Value = 0.5708267950731186
$ws.Range("K2").Value = 0.1195666986043687

$ws.Range("A3").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7f84bc13a220>),
                (''model'',
                 LogisticRegression(l1_ratio=0.95, max_iter=1000,
                                    penalty=''elasticnet'', random_state=42,
                                    solver=''saga''))])'
$ws.Range("B3").Value = 0.7058823529411764
$ws.Range("C3").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7f84bc13a280>, ''scaler'': MinMaxScaler(), ''model__solver'': ''saga'', ''model__penalty'': ''elasticnet'', ''model__l1_ratio'': 0.95, ''model__class_weight'': None}'
$ws.Range("D3").Value = 0.7368421052631579
$ws.Range("E3").Value = '[1 1 0 1 0 0 1 0 1 1 1 0]'
$ws.Range("F3").Value = '[1 1 1 1 1 1 1 1 1 1 1 1]'
$ws.Range("G3").Value = 69
$ws.Range("H3").Value = 0.7206601951574058
$ws.Range("I3").Value = 0.01863472491005711
$ws.Range("J3").Value = 0.6459373316634346
$ws.Range("K3").Value = 0.1097797140122413

$ws.Range("A4").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),
                (''model'',
                 LogisticRegression(l1_ratio=0.7, max_iter=1000,
                                    penalty=''elasticnet'', random_state=42,
                                    solver=''saga''))])'
$ws.Range("B4").Value = 0.7296078431372548
$ws.Range("C4").Value = '{''selector'': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), ''scaler'': MinMaxScaler(), ''model__solver'': ''saga'', ''model__penalty'': ''elasticnet'', ''model__l1_ratio'': 0.7, ''model__class_weight'': None}'
$ws.Range("D4").Value = 0.7777777777777777
$ws.Range("E4").Value = '[0 1 0 0 1 1 1 1 1 1 1 0]'
$ws.Range("F4").Value = '[1 1 1 1 1 1 0 1 1 1 1 0]'
$ws.Range("G4").Value = 23
$ws.Range("H4").Value = 0.7176654084892338
$ws.Range("I4").Value = 0.02092927649931291
$ws.Range("J4").Value = 0.6611737738024502
$ws.Range("K4").Value = 0.08846323748931456

$ws.Range("A5").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),
                (''model'',
                 LogisticRegression(l1_ratio=0.3, max_iter=1000,
                                    penalty=''elasticnet'', random_state=42,
                                    solver=''saga''))])'
$ws.Range("B5").Value = 0.7841830065359476
$ws.Range("C5").Value = '{''selector'': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), ''scaler'': MinMaxScaler(), ''model__solver'': ''saga'', ''model__penalty'': ''elasticnet'', ''model__l1_ratio'': 0.3, ''model__class_weight'': None}'
$ws.Range("D5").Value = 0.5882352941176471
$ws.Range("E5").Value = '[0 1 1 0 0 1 0 0 0 0 1 1]'
$ws.Range("F5").Value = '[1 1 1 1 1 1 1 1 1 1 1 1]'
$ws.Range("G5").Value = 99
$ws.Range("H5").Value = 0.7314862579019984
$ws.Range("I5").Value = 0.01586547553482015
$ws.Range("J5").Value = 0.7022387103637104
$ws.Range("K5").Value = 0.08565157848675527
